$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(167).Insert()

$ws.Range("A167").Value = 6
$ws.Range("B167").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C167").Value = "Metropolitana"
$ws.Range("D167").Value = 44694
$ws.Range("E167").Value = 13
$ws.Range("F167").Value = 100112029
$ws.Range("G167").Value = "Orégano"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 32
$ws.Range("K167").Value = 15000
$ws.Range("L167").Value = 16000
$ws.Range("M167").Value = 15469
$ws.Range("N167").Value = '$/docena de atados'
$ws.Range("O167").Value = "Región Metropolitana"
$ws.Range("P167").Value = 5156
$ws.Range("Q167").Value = 3
$ws.Range("R167").Value = "Hortaliza"
